# "carga db y modificaciones en raiz"
# - Renames two existing "playa" shared strings (space -> hyphen)
# - Adds a new row (idplaya=3, playa="Playa-Compactacion")
# - Widens column B to fit the longer text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing values (reuses the same shared-string slots in place)
$ws.Range("B2").Value = "Playa-Sarmiento"
$ws.Range("B3").Value = "Playa-Rio Cuarto"

# New 4th row
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Playa-Compactacion"

# Widen column B (stored xml width 19 <=> ColumnWidth ~18.1666667)
$ws.Columns.Item(2).ColumnWidth = 18.1666666666666667
